$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.958.80'
$ws.Range("E2").Value = '  +2.25%  '
$ws.Range("D3").Value = '2.050.41'
$ws.Range("E3").Value = '  +1.29%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '229.03'
$ws.Range("E5").Value = '  +1.13%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.616'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.69'
$ws.Range("E7").Value = '  +6.74%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("E9").Value = '  +1.98%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0808'
$ws.Range("E10").Value = '  +2.98%  '
$ws.Range("E11").Value = '  +1.08%  '
$ws.Range("D12").Value = '2.353.67'
$ws.Range("E12").Value = '  +1.66%  '
$ws.Range("E13").Value = '  +3.04%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.81'
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.751'
$ws.Range("E15").Value = '  +1.19%  '
$ws.Range("B16").Value = 'Polkadot'
$ws.Range("C16").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.28'
$ws.Range("E16").Value = '  +1.99%  '
$ws.Range("D17").Value = '2.052.99'
$ws.Range("E17").Value = '  +1.63%  '
$ws.Range("D18").Value = '37.887.16'
$ws.Range("E18").Value = '  +2.20%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.28'
$ws.Range("E19").Value = '  -3.13%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.59'
$ws.Range("E20").Value = '  +1.18%  '
$ws.Range("E21").Value = '  +2.38%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '224.49'
$ws.Range("E22").Value = '  +0.50%  '
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("E24").Value = '  -0.37%  '
$ws.Range("E25").Value = '  +2.27%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.29'
$ws.Range("E26").Value = '  +0.94%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '166.24'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.133'
$ws.Range("E28").Value = '  +4.00%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.00'
$ws.Range("E29").Value = '  +1.70%  '
$ws.Range("E30").Value = '  +1.07%  '
$ws.Range("E31").Value = '  +1.45%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.51'
$ws.Range("E32").Value = '  +0.27%  '
$ws.Range("B33").Value = 'WEMIXToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.06'
$ws.Range("E33").Value = '  +10.61%  '
$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.58'
$ws.Range("E34").Value = '  +2.80%  '
$ws.Range("E35").Value = '  +0.11%  '
$ws.Range("E36").Value = '  -0.76%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.09'
$ws.Range("E37").Value = '  +10.23%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.27'
$ws.Range("E38").Value = '  +4.81%  '
$ws.Range("E39").Value = '  -0.01%  '
$ws.Range("D40").Value = '1.489.50'
$ws.Range("E40").Value = '  +1.47%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0216'
$ws.Range("E41").Value = '  +1.12%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '96.87'
$ws.Range("E42").Value = '  +1.47%  '
$ws.Range("E43").Value = '  +2.87%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.50'
$ws.Range("E44").Value = '  +0.71%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0921'
$ws.Range("E45").Value = '  +1.33%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.15'
$ws.Range("E46").Value = '  +0.89%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.13'
$ws.Range("E47").Value = '  +12.72%  '
$ws.Range("E48").Value = '  +0.54%  '
$ws.Range("E49").Value = '  +1.09%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.09'
$ws.Range("E50").Value = '  -2.29%  '
$ws.Range("D51").Value = '2.243.43'
$ws.Range("E51").Value = '  +1.51%  '
